$p = $ppt.ActivePresentation

# --- Slide 12: title "WE Mag Points for Port wines By Country" -> "WE Mag Points By Country"
$s12 = $p.Slides.Item(12)
$tr12 = $s12.Shapes.Item(1).TextFrame.TextRange
$tr12.Paragraphs(1,1).Runs(1,1).Text = "WE Mag Points By Country"

# --- Slide 13: title paragraphs "Alcohol Content" / "Of " / "Port wine" / "By Country"
#     -> "Alcohol Content" / "By Country"  (drop the middle two paragraphs)
$s13 = $p.Slides.Item(13)
$tr13 = $s13.Shapes.Item(1).TextFrame.TextRange
$tr13.Paragraphs(2,2).Delete()
$tr13.Paragraphs(2,2).Delete()

# --- Slide 14: title paragraphs "Wine Price of " / "Port Wine By Country"
#     -> single paragraph "Wine Price By Country"
$s14 = $p.Slides.Item(14)
$tr14 = $s14.Shapes.Item(1).TextFrame.TextRange
$tr14.Paragraphs(1,1).Delete()
$tr14.Paragraphs(1,1).Runs(1,1).Text = "Wine Price By Country"
